$wb = $excel.ActiveWorkbook

$wsMeta = $wb.Worksheets.Item("Metadata")
$wsCodes = $wb.Worksheets.Item("Include from Ferlab.bio CodeS")

# --- Update the "Status" metadata value: active -> draft ---
$wsMeta.Range("B6").Value = "draft"

# --- Update the "Date" metadata value to the new publication timestamp ---
$wsMeta.Range("B8").Value = "2023-08-01T16:12:28+00:00"

# --- Ensure the header row and body rows keep their vertical-top / wrap-text
#     alignment explicitly applied (applyAlignment) on every styled cell ---
$wsMeta.Range("A1:B1").WrapText = $true
$wsMeta.Range("A1:B1").VerticalAlignment = -4160

$wsMeta.Range("A2:B14").WrapText = $true
$wsMeta.Range("A2:B14").VerticalAlignment = -4160

$wsCodes.Range("A1").WrapText = $true
$wsCodes.Range("A1").VerticalAlignment = -4160

$wsCodes.Range("A2").WrapText = $true
$wsCodes.Range("A2").VerticalAlignment = -4160

$wsCodes.Range("A3:B4").WrapText = $true
$wsCodes.Range("A3:B4").VerticalAlignment = -4160
